$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one Espárragos price-list row per record. This week's
# update inserts two new daily records at the top of the data block
# (new rows 8 and 9), which pushes every existing record (old rows 8-29)
# down by two rows (new rows 10-31). Inserting whole rows reproduces that
# shift (and the resulting dimension A1:R29 -> A1:R31) in one step.
$ws.Rows("8:9").Insert()

# New row 8: Primera, Sin especificar, week of 44487
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 44487
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 300000000
$ws.Range("G8").Value = "Espárragos"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = 1500
$ws.Range("N8").Value = "$/kilo"
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 1500
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"

# New row 9: Segunda, Sin especificar, week of 44487
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 44487
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = 300000000
$ws.Range("G9").Value = "Espárragos"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 1200
$ws.Range("L9").Value = 1200
$ws.Range("M9").Value = 1200
$ws.Range("N9").Value = "$/kilo"
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 1200
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"
